# Word COM-interop script implementing the diff:
#  - Font rename: TimesNewToman -> Times New Roman (everywhere)
#  - Title / author / email text changes
#  - Body paragraph content swapped from "quantum computing" essay to
#    "chemistry" essay, including two places where extra sentences are
#    appended at the end of existing paragraphs
#  - A new empty paragraph added at the very end of the document body

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1. Global font fix: TimesNewToman -> Times New Roman for every run in
#    the document (but do NOT touch paragraph-mark run properties).
# ---------------------------------------------------------------------
$fullRange = $d.Range(0, $d.Content.End)
$fullRange.Font.Name = "Times New Roman"

# ---------------------------------------------------------------------
# Helper: simple literal find & replace across the whole document.
# ---------------------------------------------------------------------
function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# ---------------------------------------------------------------------
# 2. Title
# ---------------------------------------------------------------------
Replace-Text "Quantum Computing: A Revolution in Information Processing" "The Realm of Chemistry: Unveiling the Secrets of Molecular Interactions"

# ---------------------------------------------------------------------
# 3. Author line: "Dr" + "." + " Ryan Edwards" (3 runs) -> "Professor Edward Wilson" (1 run)
# ---------------------------------------------------------------------
Replace-Text "Dr. Ryan Edwards" "Professor Edward Wilson"

# ---------------------------------------------------------------------
# 4. Email line
# ---------------------------------------------------------------------
Replace-Text "ryan" "wilson"
Replace-Text "edwards@quantumcomputing" "edward@schoolmail"

# ---------------------------------------------------------------------
# 5. Body paragraph 1 (the long multi-sentence paragraph)
# ---------------------------------------------------------------------
Replace-Text "The realm of quantum computing holds boundless potential for revolutionizing our comprehension of the universe and expanding the frontiers of technological advancement" "In the vast expanse of scientific exploration, chemistry stands as a pivotal discipline, delving into the intricate world of matter and its transformations"

Replace-Text " This paradigm-shifting technology harnesses the perplexing intricacies of quantum mechanics to perform computations beyond the capabilities of traditional computers" " Chemistry is the study of the composition, structure, properties, and behavior of matter"

Replace-Text " Unlike classical bits, which can only exist in states of 0 or 1, quantum bits, or qubits, can occupy superpositions of both states simultaneously" " It encompasses the interactions between atoms and molecules, unlocking the secrets of chemical bonding, reactivity, and energy"

Replace-Text " This remarkable property, known as superposition, enables quantum computers to solve problems exponentially faster than their classical counterparts" " Through chemistry, we unravel the fundamental principles governing the behavior of substances, enabling us to comprehend the intricate symphony of chemical reactions that shape our universe"

Replace-Text "The allure of quantum computing extends far beyond mere speed" "Chemistry plays a pivotal role in our daily lives, touching every aspect of modern existence"

Replace-Text " Quantum entanglement, another cornerstone of this technology, allows qubits to become interconnected in such a way that the state of one qubit instantaneously influences the state of the others, irrespective of the intervening distance" " From the food we consume to the materials we use in construction, clothing, and medicine, chemistry is the driving force behind the intricate web of chemical processes that sustain our world"

Replace-Text " This phenomenon, akin to an omniscient network, could herald a new era of ultrafast and ultrasecure communication and cryptography" " It empowers us to harness the power of chemical reactions for diverse applications, ranging from generating energy to synthesizing life-saving pharmaceuticals"

Replace-Text "The potential applications of quantum computing are vast and multifaceted" "The pursuit of chemistry is an endeavor of exploration and discovery, inviting us to unravel the mysteries of the molecular realm"

Replace-Text " From simulating complex molecular structures and accelerating drug discovery to optimizing financial models with unprecedented accuracy and harnessing the power of artificial intelligence in ways never before imagined, this technology promises to transform industries and sectors across the spectrum" " It challenges us to delve into the fundamental laws governing the interactions of matter, pushing the boundaries of our knowledge and understanding"

# Two new sentences appended at the very end of this paragraph (after the
# final "." that used to end it).
$p1 = $d.Paragraphs(5)
$p1end = $p1.Range.End - 1
$insertRange = $d.Range($p1end, $p1end)
$insertRange.InsertAfter(". Chemistry provides a gateway to comprehending the intricate choreography of atoms and molecules, revealing the hidden patterns and mechanisms that orchestrate the chemical transformations that shape our world.")
$insertRange.Font.Name = "Times New Roman"
$insertRange.Font.Size = 12
$insertRange.Font.Color = 0

# ---------------------------------------------------------------------
# 6. Summary heading paragraph text stays "Summary" - unchanged.
# ---------------------------------------------------------------------

# ---------------------------------------------------------------------
# 7. Summary body paragraph
# ---------------------------------------------------------------------
Replace-Text "Quantum computing represents a paradigm shift in information processing, leveraging the enigmatic principles of quantum mechanics to perform computations that elude classical computers" "Chemistry is the scientific study of matter and its interactions"

Replace-Text " By harnessing the power of superposition and entanglement, quantum computers hold the promise of exponential speedup, enhanced security, and transformative applications across diverse fields" " It encompasses the composition, structure, properties, and behavior of matter, revealing the intricate world of chemical bonding, reactivity, and energy"

Replace-Text " This innovative technology has the potential to reshape industries, redefine scientific frontiers, and usher in a new era of technological marvels and human understanding" " Chemistry is a fundamental discipline that plays a pivotal role in our daily lives, impacting various aspects from food production to material synthesis and medicine development"

# Two new sentences appended at the very end of the Summary paragraph.
$p2 = $d.Paragraphs(7)
$p2end = $p2.Range.End - 1
$insertRange2 = $d.Range($p2end, $p2end)
$insertRange2.InsertAfter(". Its pursuit invites exploration and discovery, inviting us to unravel the mysteries of the molecular realm and comprehend the fundamental laws governing the interactions of matter.")
$insertRange2.Font.Name = "Times New Roman"
$insertRange2.Font.Color = 0

# ---------------------------------------------------------------------
# 8. A new, empty paragraph is added at the very end of the document body.
# ---------------------------------------------------------------------
$endRange = $d.Range($d.Content.End - 1, $d.Content.End - 1)
$endRange.InsertParagraphAfter()
